$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.504.25"
$ws.Range("E2").Value = "  -2.83%  "
$ws.Range("D3").Value = "2.432.87"
$ws.Range("E3").Value = "  +8.14%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "295.44"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.72"
$ws.Range("E6").Value = "  -5.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.566"
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.504"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.75"
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.07"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "2.804.85"
$ws.Range("E14").Value = "  +8.10%  "
$ws.Range("D15").Value = "2.436.30"
$ws.Range("E15").Value = "  +8.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.842"
$ws.Range("E16").Value = "  +7.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.14"
$ws.Range("E17").Value = "  +4.93%  "
$ws.Range("D18").Value = "45.416.55"
$ws.Range("E18").Value = "  -3.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.38"
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").Value = "0.0₃0943"
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.19"
$ws.Range("E21").Value = "  +7.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.30"
$ws.Range("E22").Value = "  +3.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.17"
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.80"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("E25").Value = "  +5.43%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.22"
$ws.Range("E28").Value = "  -8.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.68"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.83"
$ws.Range("E30").Value = "  +19.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.30"
$ws.Range("E31").Value = "  +6.49%  "
$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.74"
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "148.89"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.50"
$ws.Range("E34").Value = "  +2.88%  "
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.99"
$ws.Range("E36").Value = "  +18.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.114"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.116"
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.68"
$ws.Range("E39").Value = "  -8.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.78"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.011.43"
$ws.Range("E42").Value = "  +13.94%  "
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.22"
$ws.Range("E43").Value = "  +2.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "89.00"
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.23"
$ws.Range("E46").Value = "  +27.06%  "
$ws.Range("E47").Value = "  -9.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.61"
$ws.Range("E48").Value = "  +10.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.17"
$ws.Range("E49").Value = "  +8.33%  "
$ws.Range("D50").Value = "2.675.84"
$ws.Range("E50").Value = "  +8.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.183"
$ws.Range("E51").Value = "  -0.83%  "
